# VM run GDDrec for chart_2
# Adds a new results worksheet "GDDrec_chart_2_b_20230615_11275" (a filled-in
# copy of the "Template" sheet for the chart_2_b / GDDrec run) and a matching
# scatter chart plotting its DDMIN table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "GDDrec_chart_2_b_20230615_11275"

# ---------------------------------------------------------------------------
# 2. Header block.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "ALGORITHM"
$ws.Range("B1").Value = "GDDrec"
$ws.Range("C1").Value = "TEST CASE"
$ws.Range("D1").Value = "chart_2_b"

$ws.Range("A3").Value = "OPTIONS"

# ---------------------------------------------------------------------------
# 3. OPTIONS block.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "m_modulePath"
$ws.Range("B4").Value = "/home/lukasbosshart/workspace/defects4j/bugs/chart_2_b"

$ws.Range("A5").Value = "m_sourceFolderPath"
$ws.Range("B5").Value = "source/"

$ws.Range("A6").Value = "m_unitTestFolderPath"
$ws.Range("B6").Value = "tests/"

$ws.Range("A7").Value = "m_unitTestMethod"
$ws.Range("B7").Value = "org.jfree.data.general.junit.DatasetUtilitiesTests#testBug2849731_2"

$ws.Range("A8").Value = "m_expectedResult"
$ws.Range("B8").Value = "java.lang.NullPointerException"

$ws.Range("A9").Value = "m_compilationType"
$ws.Range("B9").Value = "IN_MEMORY"

$ws.Range("A10").Value = "m_logLevel"
$ws.Range("B10").Value = "INFO"

$ws.Range("A11").Value = "m_logCompilationErrors"
$ws.Range("B11").Value = $false

$ws.Range("A12").Value = "m_logRuntimeErrors"
$ws.Range("B12").Value = $false

$ws.Range("A13").Value = "m_multipleRuns"
$ws.Range("B13").Value = $false

$ws.Range("A14").Value = "m_numberOfThreads"
$ws.Range("B14").Value = 16

$ws.Range("A15").Value = "m_preSliceCode"
$ws.Range("B15").Value = $false

$ws.Range("A16").Value = "m_graphAlgorithmFragmentLimit"
$ws.Range("B16").Value = 0

$ws.Range("A17").Value = "m_graphAlgorithmEscalatingFragmentLimit"
$ws.Range("B17").Value = $false

# ---------------------------------------------------------------------------
# 4. MEASUREMENTS block.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "MEASUREMENTS"

$ws.Range("A20").Value = "Input Size (Fragments)"
$ws.Range("B20").Value = 407987
$ws.Range("C20").Value = "Number of DDmin runs"
$ws.Range("D20").Formula = "=COUNTA(A32:A33)"

$ws.Range("A21").Value = "Input Size (bytes)"
$ws.Range("B21").Value = 6838600
$ws.Range("C21").Value = "Average CC per run"
$ws.Range("D21").Formula = "=B27/D20"

$ws.Range("A22").Value = "Input Size (LoC)"
$ws.Range("B22").Value = 227306
$ws.Range("C22").Value = "Total Reduction Factor"
$ws.Range("D22").Formula = "=(B20-B23)/B20"
$ws.Range("D22").NumberFormat = "0.00%"

$ws.Range("A23").Value = "Output Size (Fragments)"
$ws.Range("B23").Value = 393931
$ws.Range("C23").Value = "Average Reduction per run"
$ws.Range("D23").Formula = "=D22/D20"
$ws.Range("D23").NumberFormat = "0.00%"

$ws.Range("A24").Value = "Output Size (bytes)"
$ws.Range("B24").Value = 3026209
$ws.Range("C24").Value = "Execution Time"
$ws.Range("D24").Formula = '=TEXT(B26/86400000,"hh:mm:ss.000")'

$ws.Range("A25").Value = "Output Size (LoC)"
$ws.Range("B25").Value = 108815
$ws.Range("C25").Value = "Average run size"
$ws.Range("D25").Formula = "=AVERAGE(B32:B33)"

$ws.Range("A26").Value = "Execution Time (ms)"
$ws.Range("B26").Value = 14409572

$ws.Range("A27").Value = "Compiler Calls"
$ws.Range("B27").Value = 11716

$ws.Range("A28").Value = "Failed Test Runs"
$ws.Range("B28").Value = 87

# ---------------------------------------------------------------------------
# 5. DDMIN table.
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = "DDMIN"

$ws.Range("A31").Value = "Identifier"
$ws.Range("B31").Value = "Active Number"
$ws.Range("C31").Value = "Active Result"
$ws.Range("D31").Value = "Active Reduction"
$ws.Range("E31").Value = "Total Number"
$ws.Range("F31").Value = "Total Result"
$ws.Range("G31").Value = "Total Reduction"
$ws.Range("H31").Value = "Compiler Calls"
$ws.Range("I31").Value = "Failed Runs"
$ws.Range("J31").Value = "Overall CC"
$ws.Range("K31").Value = "Overall FR"
$ws.Range("L31").Value = "Time (ms)"

$ws.Range("A32").Value = "0-0"
$ws.Range("B32").Value = 654
$ws.Range("C32").Value = 611
$ws.Range("D32").Formula = "=B32-C32"
$ws.Range("E32").Value = 407987
$ws.Range("F32").Value = 393931
$ws.Range("G32").Formula = "=E32-F32"
$ws.Range("H32").Value = 11716
$ws.Range("I32").Value = 87
$ws.Range("J32").Formula = "=SUM($H$32:H32)"
$ws.Range("K32").Formula = "=SUM($I$32:I32)"
$ws.Range("L32").Value = 13754400
$ws.Range("M32").Value = 14408143

$ws.Range("A33").Value = "0-1"
$ws.Range("B33").Value = 10
$ws.Range("C33").Value = 10
$ws.Range("D33").Formula = "=B33-C33"
$ws.Range("E33").Value = 393931
$ws.Range("F33").Value = 393931
$ws.Range("G33").Formula = "=E33-F33"
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Formula = "=SUM($H$32:H33)"
$ws.Range("K33").Formula = "=SUM($I$32:I33)"
$ws.Range("L33").Value = 103
$ws.Range("M33").Value = 14408246
$ws.Rows.Item(33).RowHeight = 15

# ---------------------------------------------------------------------------
# 6. Formatting: bold + vertically-centered headers, bold labels, column
#    widths copied from the Template sheet this was derived from.
# ---------------------------------------------------------------------------
$boldCentered = $ws.Range("A1,C1,A19,A30,A31:L31")
$boldCentered.Font.Bold = $true
$boldCentered.VerticalAlignment = -4108

$centered = $ws.Range("A2,A3,A16:A18")
$centered.VerticalAlignment = -4108

$labelRange = $ws.Range("A4:A15,A20:A28")
$labelRange.VerticalAlignment = -4108

$ws.Columns.Item(1).ColumnWidth = 39.86
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(3).ColumnWidth = 21.57
$ws.Columns.Item(4).ColumnWidth = 16.14
$ws.Columns.Item(5).ColumnWidth = 18.86
$ws.Columns.Item(6).ColumnWidth = 13.71
$ws.Columns.Item(7).ColumnWidth = 14.86
$ws.Columns.Item(8).ColumnWidth = 13.71

$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 7. Scatter chart: "# of fragments after # of compiler calls", plotting the
#    new sheet's DDMIN table (mirrors the chart already sitting on Template).
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Add(0, 0, 300, 200)
$co.Name = "Chart 1"
$chart = $co.Chart
$chart.ChartType = 74

$chart.HasTitle = $true
$chart.ChartTitle.Text = "# of fragments after # of compiler calls"

$chart.SeriesCollection().NewSeries()
$ser = $chart.SeriesCollection().Item(1)
$ser.Name = "=Template!`$F`$31"
$ser.XValues = "=GDDrec_chart_2_b_20230615_11275!`$J`$32:`$J`$34"
$ser.Values = "=GDDrec_chart_2_b_20230615_11275!`$F`$32:`$F`$34"

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "# compiler calls"
$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "# fragments"

$chart.HasLegend = $false

$excel.CalculateFull()
